$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Table 1 (header table with tblInd 34 -> 26): only indentation / left
# cell-margin values change; borders stay as-is.
# ---------------------------------------------------------------------------
$t1 = $d.Tables(1)

# w:tblInd 34 -> 26 (dxa) == 1.7pt -> 1.3pt
$t1.Rows.LeftIndent = 1.3

# w:tblCellMar left 19 -> 11 (dxa) == 0.95pt -> 0.55pt
$t1.LeftPadding = 0.55

# Each of the 4 cells (2 rows x 2 cols) has its own left tcMar override,
# 19 -> 11 dxa == 0.95pt -> 0.55pt
for ($r = 1; $r -le $t1.Rows.Count; $r++) {
    for ($c = 1; $c -le $t1.Columns.Count; $c++) {
        $t1.Cell($r, $c).LeftPadding = 0.55
    }
}

# ---------------------------------------------------------------------------
# Table 2 (content table with tblInd 32 -> 30): indentation, left cell
# margins, and every border width (sz 2 -> sz 6, i.e. 0.25pt -> 0.75pt).
# ---------------------------------------------------------------------------
$t2 = $d.Tables(2)

# w:tblInd 32 -> 30 (dxa) == 1.6pt -> 1.5pt
$t2.Rows.LeftIndent = 1.5

# w:tblCellMar left 25 -> 17 (dxa) == 1.25pt -> 0.85pt
$t2.LeftPadding = 0.85

# Table-level borders: only top / left / bottom / insideH exist; update each
# individually (by WdBorderType index) so no new right/insideV edge gets
# fabricated. wdLineWidth075pt == 3 (sz=6, i.e. 0.75pt).
$t2.Borders.Item(-1).LineWidth = 3   # wdBorderTop
$t2.Borders.Item(-2).LineWidth = 3   # wdBorderLeft
$t2.Borders.Item(-3).LineWidth = 3   # wdBorderBottom
$t2.Borders.Item(-5).LineWidth = 3   # wdBorderHorizontal (inside)

# Row 1, Cell 1 (the narrow "label" cell): same 4 edges as the table, plus
# its own left tcMar override 25 -> 17 dxa == 1.25pt -> 0.85pt.
$c1 = $t2.Cell(1, 1)
$c1.LeftPadding = 0.85
$c1.Borders.Item(-1).LineWidth = 3   # wdBorderTop
$c1.Borders.Item(-2).LineWidth = 3   # wdBorderLeft
$c1.Borders.Item(-3).LineWidth = 3   # wdBorderBottom
$c1.Borders.Item(-5).LineWidth = 3   # wdBorderHorizontal (inside)

# Row 1, Cell 2 (the wide "content" cell): all 6 edges present, plus its own
# left tcMar override 25 -> 17 dxa == 1.25pt -> 0.85pt.
$c2 = $t2.Cell(1, 2)
$c2.LeftPadding = 0.85
$c2.Borders.Item(-1).LineWidth = 3   # wdBorderTop
$c2.Borders.Item(-2).LineWidth = 3   # wdBorderLeft
$c2.Borders.Item(-3).LineWidth = 3   # wdBorderBottom
$c2.Borders.Item(-4).LineWidth = 3   # wdBorderRight
$c2.Borders.Item(-5).LineWidth = 3   # wdBorderHorizontal (inside)
$c2.Borders.Item(-6).LineWidth = 3   # wdBorderVertical (inside)
